$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 (columns D-H) with new values
$ws.Range("D2").Value = 6.097791848989952
$ws.Range("E2").Value = 3.158824152795895
$ws.Range("F2").Value = 0.002098350418317629
$ws.Range("G2").Value = 0.0005304426849321864
$ws.Range("H2").Value = 0.002164357622956945

$ws.Range("D3").Value = 1.327958352405515
$ws.Range("E3").Value = 1.127961348663399
$ws.Range("F3").Value = 0.0001701609159572293
$ws.Range("G3").Value = -0.0042597355648335
$ws.Range("H3").Value = 0.004263132864411661

$ws.Range("D4").Value = 7.508419772969479
$ws.Range("E4").Value = 8.301797276737126
$ws.Range("F4").Value = 0.002257859367300163
$ws.Range("G4").Value = -0.0002473414846484445
$ws.Range("H4").Value = 0.002271366710272296

# Add new row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5.853086206137439
$ws.Range("C5").Value = 6.835972487871987
$ws.Range("D5").Value = 5.848431766679135
$ws.Range("E5").Value = 6.832672373997851
$ws.Range("F5").Value = -0.004654439458303905
$ws.Range("G5").Value = -0.003300113874136201
$ws.Range("H5").Value = 0.005705660194340581
